# The workbook originally had an extra "HC2" column (column I) on the
# "Results" sheet that has been removed from the published version. Delete
# that whole column; Excel will shift everything to its right one column to
# the left and the now-unused "HC2" shared string will drop out of the
# shared strings table automatically when the file is saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")
$ws.Range("I1").EntireColumn.Delete()

# The last time the workbook was saved, the "Results" sheet was the active
# tab (with cell K2 selected), rather than "Synthetic".
$ws.Activate()
$ws.Range("K2").Select()
